{"js": "// New replacement variable: swap the hard-coded company name literal\n// for the $CURRENTUSERCOMPANY$ merge/replacement token, the same way\n// the other signature fields (name, title, department, ...) already\n// use $CURRENTUSER...$ tokens.\n//\n// We look the literal text up (rather than hard-coding a paragraph/run\n// index) so the edit is resilient to unrelated layout differences, and\n// we insert the new text right after the matched range and then delete\n// the old range \u2014 this keeps the insertion anchored to the run that\n// carries the \"Super Duper Inc.\" text (bold, dark-grey signature\n// heading) instead of bleeding into the formatting of a neighboring\n// run.\n\nconst searchText = \"Super Duper Inc.\";\nconst replacement = \"$CURRENTUSERCOMPANY$\";\n\nconst body = context.document.body;\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const found = results.items[i];\n  found.insertText(replacement, Word.InsertLocation.after);\n  found.delete();\n}\nawait context.sync();\n", "ps1": "# New replacement variable: swap the hard-coded company name literal\n# for the $CURRENTUSERCOMPANY$ merge/replacement token, the same way\n# the other signature fields (name, title, department, ...) already\n# use $CURRENTUSER...$ tokens.\n\n$d = $word.ActiveDocument\n\n$searchText = \"Super Duper Inc.\"\n$replacement = \"`$CURRENTUSERCOMPANY`$\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $searchText\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop - one match at a time, we re-scan from the top each loop\n\nwhile ($find.Execute()) {\n    # $rng tracks the found span; remember its bounds before InsertAfter\n    # expands $rng itself to also cover the newly inserted text.\n    $origStart = $rng.Start\n    $origEnd = $rng.End\n\n    $rng.InsertAfter($replacement)\n\n    # Delete only the original matched text (by position), leaving the\n    # freshly inserted replacement - and the formatting of the run it\n    # belongs to - untouched.\n    $oldRange = $d.Range($origStart, $origEnd)\n    $oldRange.Delete()\n\n    # Re-scan from the very beginning in case other occurrences exist.\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Forward = $true\n    $find.Wrap = 0\n}\n"}
